$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C (old "Quantidade" shifts to D,
# old "Quantidade Vendida" shifts to E).
$ws.Columns.Item(3).Insert()

# New column header - copy formatting from the (shifted) "Quantidade" header
# cell, then set the text.
$ws.Cells.Item(1, 4).Copy()
$ws.Cells.Item(1, 3).PasteSpecial(-4122)
$ws.Cells.Item(1, 3).Value2 = "Preço de Venda"

# Fill new column C with Preço de Venda = Valor Unitário (col B) * 1.2
for ($r = 2; $r -le 121; $r++) {
    $unitValue = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($r, 3).Value = $unitValue * 1.2
}
